$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Score")

# Restore the workbook window to its normal (non-maximized) size/position.
$win = $excel.ActiveWindow
$win.WindowState = -4143
$win.Left = 3000
$win.Top = 3000
$win.Width = 17280
$win.Height = 8928

$data = @{
  2  = @(8, 0, 258, 44, 0, 32, 28, 36, 15, 10, 15)
  3  = @(8, 0, 266, 44, 0, 35, 37, 24, 16, 31, 36)
  4  = @(8, 0, 300, 44, 0, 17, 18, 29, 34, 12, 18)
  5  = @(6, 3, 300, 44, 0, 24, 31, 32, 15, 33, 16)
  6  = @(7, 47, 300, 44, 0, 34, 25, 30, 13, 10, 15)
  7  = @(7, 43, 300, 44, 0, 16, 11, 33, 35, 26, 11)
  8  = @(5, 49, 300, 44, 0, 10, 36, 22, 21, 25, 17)
  9  = @(6, 42, 300, 44, 0, 39, 27, 27, 33, 21, 19)
  10 = @(6, 45, 300, 44, 0, 24, 27, 38, 38, 12, 11)
  11 = @(6, 45, 300, 44, 0, 13, 31, 34, 17, 12, 34)
  12 = @(6, 43, 300, 44, 0, 37, 27, 23, 23, 24, 31)
  13 = @(8, 0, 300, 44, 0, 20, 37, 26, 25, 29, 29)
  14 = @(4, 56, 300, 44, 0, 36, 31, 21, 38, 10, 18)
  15 = @(5, 58, 300, 44, 0, 25, 29, 31, 14, 15, 20)
  16 = @(6, 42, 300, 44, 0, 38, 12, 21, 37, 39, 28)
  17 = @(5, 38, 300, 44, 0, 32, 29, 27, 15, 26, 28)
  18 = @(8, 0, 300, 44, 0, 15, 14, 29, 39, 12, 11)
  19 = @(6, 15, 300, 44, 0, 26, 25, 13, 25, 27, 39)
  20 = @(7, 23, 300, 44, 0, 25, 13, 29, 32, 26, 34)
  21 = @(5, 33, 300, 44, 0, 21, 32, 34, 28, 34, 17)
}

foreach ($row in 2..21) {
  $values = $data[$row]
  $ws.Range("C$row").Value = $values[0]
  $ws.Range("D$row").Value = $values[1]
  $ws.Range("E$row").Value = $values[2]
  $ws.Range("F$row").Value = $values[3]
  $ws.Range("G$row").Value = $values[4]
  $ws.Range("H$row").Value = $values[5]
  $ws.Range("I$row").Value = $values[6]
  $ws.Range("J$row").Value = $values[7]
  $ws.Range("K$row").Value = $values[8]
  $ws.Range("L$row").Value = $values[9]
  $ws.Range("M$row").Value = $values[10]
  $ws.Range("N$row").Formula = "=SUM(H$row`:M$row)"
}

$ws.Range("G5").Select()
